$d = $word.ActiveDocument

# Delete all paragraphs after the first (bookmark) paragraph, leaving just
# the empty paragraph that carries the _GoBack bookmark.
$count = $d.Paragraphs.Count
if ($count -gt 1) {
    $start = $d.Paragraphs(2).Range.Start
    $end = $d.Paragraphs($count).Range.End
    $r = $d.Range($start, $end)
    $r.Delete()
}

# Mark the Normal style's font as superscript (adds
# <w:vertAlign w:val="superscript"/> to its run properties) so new
# citation markers inherit superscript formatting.
$normal = $d.Styles("Normal")
$normal.Font.Superscript = $true

Write-Output ("ParaCount=" + $d.Paragraphs.Count)
Write-Output ("NormalSuperscript=" + $normal.Font.Superscript)
